$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 1: A1/B1 become #N/A error literals (was inline "NA" strings)
$ws.Range("A1").Value = "#N/A"
$ws.Range("B1").Value = "#N/A"

# 2) Insert 13 new rows for the "Klebsiella pneumoniae complex" block
#    right before the existing "Meyerozyma guilliermondii complex" block (old row 314)
$ws.Range("A314:A326").EntireRow.Insert()

# 3) Populate the newly inserted rows with the Klebsiella pneumoniae complex data
$ws.Cells.Item(314, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(314, 2).Value = "B_KLBSL_AFRC"
$ws.Cells.Item(314, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(314, 4).Value = "Klebsiella africana"
$ws.Cells.Item(315, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(315, 2).Value = "B_KLBSL_PNMN"
$ws.Cells.Item(315, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(315, 4).Value = "Klebsiella pneumoniae"
$ws.Cells.Item(316, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(316, 2).Value = "B_KLBSL_PNMN_OZAN"
$ws.Cells.Item(316, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(316, 4).Value = "Klebsiella pneumoniae ozaenae"
$ws.Cells.Item(317, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(317, 2).Value = "B_KLBSL_PNMN_PNMN"
$ws.Cells.Item(317, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(317, 4).Value = "Klebsiella pneumoniae pneumoniae"
$ws.Cells.Item(318, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(318, 2).Value = "B_KLBSL_PNMN_RHNS"
$ws.Cells.Item(318, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(318, 4).Value = "Klebsiella pneumoniae rhinoscleromatis"
$ws.Cells.Item(319, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(319, 2).Value = "B_KLBSL_QSPN"
$ws.Cells.Item(319, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(319, 4).Value = "Klebsiella quasipneumoniae"
$ws.Cells.Item(320, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(320, 2).Value = "B_KLBSL_QSPN_QSPN"
$ws.Cells.Item(320, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(320, 4).Value = "Klebsiella quasipneumoniae quasipneumoniae"
$ws.Cells.Item(321, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(321, 2).Value = "B_KLBSL_QSPN_SMLP"
$ws.Cells.Item(321, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(321, 4).Value = "Klebsiella quasipneumoniae similipneumoniae"
$ws.Cells.Item(322, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(322, 2).Value = "B_KLBSL_QSVR"
$ws.Cells.Item(322, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(322, 4).Value = "Klebsiella quasivariicola"
$ws.Cells.Item(323, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(323, 2).Value = "B_KLBSL_VRCL"
$ws.Cells.Item(323, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(323, 4).Value = "Klebsiella variicola"
$ws.Cells.Item(324, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(324, 2).Value = "B_KLBSL_VRCL_TRPC"
$ws.Cells.Item(324, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(324, 4).Value = "Klebsiella variicola tropica"
$ws.Cells.Item(325, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(325, 2).Value = "B_KLBSL_VRCL_LNSS"
$ws.Cells.Item(325, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(325, 4).Value = "Klebsiella variicola tropicalensis"
$ws.Cells.Item(326, 1).Value = "B_KLBSL_PNMN-C"
$ws.Cells.Item(326, 2).Value = "B_KLBSL_VRCL_VRCL"
$ws.Cells.Item(326, 3).Value = "Klebsiella pneumoniae complex"
$ws.Cells.Item(326, 4).Value = "Klebsiella variicola variicola"
